$p = $ppt.ActivePresentation

# The deck currently uses the "Integral" design/theme (green palette) on its
# single Slide Master. This switches the applied design over to the built-in
# "Office Theme" colour palette (blue), matching the colour scheme that the
# companion "Office Theme" theme part in this deck already carries.
#
# PowerPoint's ColorScheme.Colors(index) slots map to the clrScheme children
# in document order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB values are packed as R + G*256 + B*65536 (the VBA RGB() layout).

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$colorScheme = $master.ColorScheme

$colorScheme.Colors(1).RGB  = 0          # dk1      #000000
$colorScheme.Colors(2).RGB  = 16777215   # lt1      #FFFFFF
$colorScheme.Colors(3).RGB  = 6968388    # dk2      #44546A
$colorScheme.Colors(4).RGB  = 15132391   # lt2      #E7E6E6
$colorScheme.Colors(5).RGB  = 13998939   # accent1  #5B9BD5
$colorScheme.Colors(6).RGB  = 3243501    # accent2  #ED7D31
$colorScheme.Colors(7).RGB  = 10855845   # accent3  #A5A5A5
$colorScheme.Colors(8).RGB  = 49407      # accent4  #FFC000
$colorScheme.Colors(9).RGB  = 12874308   # accent5  #4472C4
$colorScheme.Colors(10).RGB = 4697456    # accent6  #70AD47
$colorScheme.Colors(11).RGB = 12673797   # hlink    #0563C1
$colorScheme.Colors(12).RGB = 7491477    # folHlink #954F72
